$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "PlayerClass" header into the new column P (column 16),
# reusing the same text already used in I1 (PlayerClass).
$ws.Range("P1").Value = $ws.Range("I1").Value2

# Move the active selection to L18, matching the saved workbook state.
$ws.Range("L18").Select()
